$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column M entirely; the data that was in column N shifts left to
# become the new column M (effectively removing a column from the data).
$ws.Columns("M:M").Delete()

# Reflect the new active cell/selection after the column removal.
$ws.Range("M1").Select()
